$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top Gainers")

$ws.Range("C2").Value = 11.952
$ws.Range("D2").Value = 19.8443
$ws.Range("E2").Value = 33.4974

$ws.Range("C3").Value = 11.4928
$ws.Range("D3").Value = 20.3242
$ws.Range("E3").Value = 27.5023

$ws.Range("C4").Value = 9.4406
$ws.Range("D4").Value = 6.5143
$ws.Range("E4").Value = -8.2111

$ws.Range("C5").Value = 8.999599999999999
$ws.Range("D5").Value = 9.2033
$ws.Range("E5").Value = 23.0777

$ws.Range("C6").Value = 8.136200000000001
$ws.Range("D6").Value = 7.0998
$ws.Range("E6").Value = -21.6465

$ws.Range("C7").Value = 8.070600000000001
$ws.Range("D7").Value = 10.3552
$ws.Range("E7").Value = 26.6035

$ws.Range("C8").Value = 7.5763
$ws.Range("D8").Value = 7.8657
$ws.Range("E8").Value = 9.635400000000001

$ws.Range("B9").Value = "CHENNPETRO"
$ws.Range("C9").Value = 7.5213
$ws.Range("D9").Value = 12.7209
$ws.Range("E9").Value = 14.8852

$ws.Range("B10").Value = "SHREEJISPG"
$ws.Range("C10").Value = 7.1843
$ws.Range("D10").Value = 11.5018
$ws.Range("E10").Value = 12.6077

$ws.Range("C11").Value = 6.4905
$ws.Range("D11").Value = 12.4237
$ws.Range("E11").Value = 14.5375

$ws.Range("C12").Value = 6.4014
$ws.Range("D12").Value = 9.772399999999999
$ws.Range("E12").Value = 15.983

$ws.Range("C13").Value = 6.0491
$ws.Range("D13").Value = 19.4879
$ws.Range("E13").Value = 19.577

$ws.Range("C14").Value = 5.9417
$ws.Range("D14").Value = 6.5743
$ws.Range("E14").Value = 4.2562

$ws.Range("B16").Value = "BHEL"
$ws.Range("C16").Value = 5.3221
$ws.Range("D16").Value = 11.9122
$ws.Range("E16").Value = 8.387499999999999

$ws.Range("B17").Value = "MRPL"
$ws.Range("C17").Value = 5.2754
$ws.Range("D17").Value = 15.4979
$ws.Range("E17").Value = 26.3875

$ws.Range("B18").Value = "BAJAJHCARE"
$ws.Range("C18").Value = 5.2718
$ws.Range("D18").Value = 5.825
$ws.Range("E18").Value = -0.5249

$ws.Range("B21").Value = "VIMTALABS"
$ws.Range("C21").Value = 4.8132
$ws.Range("D21").Value = 4.9926
$ws.Range("E21").Value = -0.1345

$ws.Range("B22").Value = "SHANTIGOLD"
$ws.Range("C22").Value = 4.7784
$ws.Range("D22").Value = 12.1531
$ws.Range("E22").Value = 4.6706

$ws.Range("B23").Value = "RAMASTEEL"
$ws.Range("C23").Value = 4.7189
$ws.Range("D23").Value = 4.6138
$ws.Range("E23").Value = 6.2118

$ws.Range("B24").Value = "IIFL"
$ws.Range("C24").Value = 4.5889
$ws.Range("D24").Value = 11.6782
$ws.Range("E24").Value = 21.0445

$ws.Range("B25").Value = "DEEDEV"
$ws.Range("C25").Value = 4.4746
$ws.Range("D25").Value = -2.4759
$ws.Range("E25").Value = -3.2802

$ws.Range("C26").Value = 4.4062
$ws.Range("D26").Value = 10.2839
$ws.Range("E26").Value = 27.2541

$ws.Range("B27").Value = "MARINE"
$ws.Range("C27").Value = 4.3338
$ws.Range("D27").Value = 0.9284
$ws.Range("E27").Value = 13.2533

$ws.Range("B28").Value = "VSTIND"
$ws.Range("C28").Value = 4.3125
$ws.Range("D28").Value = 4.8193
$ws.Range("E28").Value = 4.3529

$ws.Range("B29").Value = "HIRECT"
$ws.Range("C29").Value = 4.1509
$ws.Range("D29").Value = 11.6032
$ws.Range("E29").Value = 9.8894

$ws.Range("B30").Value = "BLS"
$ws.Range("C30").Value = 3.8854
$ws.Range("D30").Value = 0.8169999999999999
$ws.Range("E30").Value = -0.4499

$ws.Range("B31").Value = "SALASAR"
$ws.Range("C31").Value = 3.8579
$ws.Range("D31").Value = 8.829800000000001
$ws.Range("E31").Value = 15.3326

$ws.Range("B32").Value = "CENTRUM"
$ws.Range("C32").Value = 3.6667
$ws.Range("D32").Value = 2.272
$ws.Range("E32").Value = 1.6038

$ws.Range("B33").Value = "SKYGOLD"
$ws.Range("C33").Value = 3.6475
$ws.Range("D33").Value = -0.9026
$ws.Range("E33").Value = 37.6258

$ws.Range("B34").Value = "POLICYBZR"
$ws.Range("C34").Value = 3.5453
$ws.Range("D34").Value = 5.8611
$ws.Range("E34").Value = 4.8472

$ws.Range("C35").Value = 3.4076
$ws.Range("D35").Value = 4.6169
$ws.Range("E35").Value = 24.8256

$ws.Range("B36").Value = "OIL"
$ws.Range("C36").Value = 3.365
$ws.Range("D36").Value = 3.6114
$ws.Range("E36").Value = 5.0387

$ws.Range("B37").Value = "MTARTECH"
$ws.Range("C37").Value = 3.3022
$ws.Range("D37").Value = 7.4153
$ws.Range("E37").Value = 31.2898

$ws.Range("B38").Value = "RSYSTEMS"
$ws.Range("C38").Value = 3.2819
$ws.Range("D38").Value = 4.4412
$ws.Range("E38").Value = 6.7864

$ws.Range("B39").Value = "RELTD"
$ws.Range("C39").Value = 3.2635
$ws.Range("D39").Value = 10.0542
$ws.Range("E39").Value = -1.4749

$ws.Range("B40").Value = "CENTUM"
$ws.Range("C40").Value = 3.2333
$ws.Range("D40").Value = 4.0951
$ws.Range("E40").Value = -1.4238

$ws.Range("B41").Value = "ABREL"
$ws.Range("C41").Value = 3.1687
$ws.Range("D41").Value = 11.4717
$ws.Range("E41").Value = 11.0098

$ws.Range("C42").Value = 3.1575
$ws.Range("D42").Value = 7.5565
$ws.Range("E42").Value = 0.1833

$ws.Range("B43").Value = "GPPL"
$ws.Range("C43").Value = 3.0661
$ws.Range("D43").Value = 6.5778
$ws.Range("E43").Value = 8.2706

$ws.Range("B44").Value = "BPCL"
$ws.Range("C44").Value = 3.0595
$ws.Range("D44").Value = 8.5641
$ws.Range("E44").Value = 5.6234

$ws.Range("B45").Value = "REFEX"
$ws.Range("C45").Value = 2.9619
$ws.Range("D45").Value = 0.385
$ws.Range("E45").Value = 2.3699

$ws.Range("B47").Value = "DBCORP"
$ws.Range("C47").Value = 2.8055
$ws.Range("D47").Value = 5.4461
$ws.Range("E47").Value = 1.3928

$ws.Range("B48").Value = "CREDITACC"
$ws.Range("C48").Value = 2.7855
$ws.Range("D48").Value = 1.4271
$ws.Range("E48").Value = 6.6213

$ws.Range("B49").Value = "MAMATA"
$ws.Range("C49").Value = 2.752
$ws.Range("D49").Value = 2.1619
$ws.Range("E49").Value = 1.363

$ws.Range("B50").Value = "MFSL"
$ws.Range("C50").Value = 2.7127
$ws.Range("D50").Value = 2.7668
$ws.Range("E50").Value = -1.0089

$ws.Range("B51").Value = "CIFL"
$ws.Range("C51").Value = 2.7043
$ws.Range("D51").Value = 2.2582
$ws.Range("E51").Value = 2.1991

$ws.Range("B52").Value = "GANESHCP"
$ws.Range("C52").Value = 2.6912
$ws.Range("D52").Value = 2.1666
$ws.Range("E52").Value = 1.7187

$ws.Range("B53").Value = "NBCC"
$ws.Range("C53").Value = 2.6818
$ws.Range("D53").Value = 5.9271
$ws.Range("E53").Value = 10.4875

$ws.Range("B54").Value = "OBEROIRLTY"
$ws.Range("C54").Value = 2.6614
$ws.Range("D54").Value = 3.4767
$ws.Range("E54").Value = 11.1813

$ws.Range("B55").Value = "SDBL"
$ws.Range("C55").Value = 2.6379
$ws.Range("D55").Value = 1.182
$ws.Range("E55").Value = 6.7761

$ws.Range("B56").Value = "GENUSPOWER"
$ws.Range("C56").Value = 2.6218
$ws.Range("D56").Value = 10.5867
$ws.Range("E56").Value = 7.3053

$ws.Range("C57").Value = 2.5969
$ws.Range("D57").Value = -5.0051
$ws.Range("E57").Value = 10.873

$ws.Range("B58").Value = "PSPPROJECT"
$ws.Range("C58").Value = 2.5111
$ws.Range("D58").Value = 16.5354
$ws.Range("E58").Value = 22.9424

$ws.Range("B59").Value = "INOXGREEN"
$ws.Range("C59").Value = 2.5061
$ws.Range("D59").Value = 10.498
$ws.Range("E59").Value = 33.8273

$ws.Range("B60").Value = "GMRAIRPORT"
$ws.Range("C60").Value = 2.5062
$ws.Range("D60").Value = 2.2532
$ws.Range("E60").Value = 9.289

$ws.Range("B61").Value = "CARYSIL"
$ws.Range("C61").Value = 2.4914
$ws.Range("D61").Value = 1.9765
$ws.Range("E61").Value = 10.8493

$ws.Range("B62").Value = "DIVISLAB"
$ws.Range("C62").Value = 2.4877
$ws.Range("D62").Value = 1.2286
$ws.Range("E62").Value = 17.3038

$ws.Range("B63").Value = "BGRENERGY"
$ws.Range("C63").Value = 2.4727
$ws.Range("D63").Value = -6.7673
$ws.Range("E63").Value = 73.5703

$ws.Range("B64").Value = "DBL"
$ws.Range("C64").Value = 2.4656
$ws.Range("D64").Value = 3.549
$ws.Range("E64").Value = 4.6224

$ws.Range("B65").Value = "ALICON"
$ws.Range("C65").Value = 2.4633
$ws.Range("D65").Value = 8.5374
$ws.Range("E65").Value = 13.9288

$ws.Range("B66").Value = "MOTILALOFS"
$ws.Range("C66").Value = 2.4509
$ws.Range("D66").Value = -0.3571
$ws.Range("E66").Value = 15.5439

$ws.Range("B67").Value = "SOLEX"
$ws.Range("C67").Value = 2.4418
$ws.Range("D67").Value = 3.9749
$ws.Range("E67").Value = "N/A"

$ws.Range("B68").Value = "CAMS"
$ws.Range("C68").Value = 2.4381
$ws.Range("D68").Value = 1.786
$ws.Range("E68").Value = 5.079

$ws.Range("B69").Value = "ACUTAAS"
$ws.Range("C69").Value = 2.4093
$ws.Range("D69").Value = 2.7743
$ws.Range("E69").Value = 31.2794

$ws.Range("B70").Value = "BLISSGVS"
$ws.Range("C70").Value = 2.4005
$ws.Range("D70").Value = 1.7581
$ws.Range("E70").Value = 2.1251

$ws.Range("B71").Value = "CEATLTD"
$ws.Range("C71").Value = 2.3463
$ws.Range("D71").Value = -0.7704
$ws.Range("E71").Value = 21.4459

$ws.Range("B72").Value = "NEULANDLAB"
$ws.Range("C72").Value = 2.3238
$ws.Range("D72").Value = -2.0139
$ws.Range("E72").Value = 7.9356

$ws.Range("B73").Value = "IFCI"
$ws.Range("C73").Value = 2.3129
$ws.Range("D73").Value = 4.8859
$ws.Range("E73").Value = 8.0296

$ws.Range("B74").Value = "ANANDRATHI"
$ws.Range("C74").Value = 2.2713
$ws.Range("D74").Value = 1.3843
$ws.Range("E74").Value = 11.6634

$ws.Range("B75").Value = "AARTIDRUGS"
$ws.Range("C75").Value = 2.2472
$ws.Range("D75").Value = 2.2472
$ws.Range("E75").Value = 3.8129

$ws.Range("B76").Value = "INDRAMEDCO"
$ws.Range("C76").Value = 2.2462
$ws.Range("D76").Value = 7.0261
$ws.Range("E76").Value = 30.9915
